$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.659.28"
$ws.Range("E2").Value = "  +4.13%  "
$ws.Range("D3").Value = "'1.797.55"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'313.46"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "'0.5308"
$ws.Range("E7").Value = "  -1.22%  "
$ws.Range("D8").Value = "'0.3774"
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("D9").Value = "'0.07542"
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("D10").Value = "'42.52"
$ws.Range("E10").Value = "  -1.14%  "
$ws.Range("D11").Value = "'1.122"
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("D12").Value = "'21.19"
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("D13").Value = "'1.001"
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("D14").Value = "'6.201"
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("D15").Value = "'7.458"
$ws.Range("E15").Value = "  +5.98%  "
$ws.Range("D16").Value = "'1.791.08"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").Value = "'90.53"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").Value = "'0.00001067"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").Value = "'0.06443"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").Value = "'0.9999"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "'17.30"
$ws.Range("E21").Value = "  +2.21%  "
$ws.Range("D22").Value = "'5.933"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "'28.659.80"
$ws.Range("E23").Value = "  +4.15%  "
$ws.Range("D24").Value = "'11.20"
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("D25").Value = "'2.095"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "'160.79"
$ws.Range("E26").Value = "  +2.82%  "
$ws.Range("D27").Value = "'20.60"
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("D28").Value = "'2.410"
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").Value = "'1.998.36"
$ws.Range("E29").Value = "  +0.70%  "
$ws.Range("D30").Value = "'123.64"
$ws.Range("E30").Value = "  +1.56%  "
$ws.Range("D31").Value = "'1.132"
$ws.Range("E31").Value = "  +3.11%  "
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").Value = "'5.740"
$ws.Range("E33").Value = "  +1.82%  "
$ws.Range("D34").Value = "'3.663"
$ws.Range("E34").Value = "  +1.11%  "
$ws.Range("E35").Value = "  +11.10%  "
$ws.Range("D36").Value = "'0.06580"
$ws.Range("E36").Value = "  +9.43%  "
$ws.Range("D37").Value = "'0.02328"
$ws.Range("E37").Value = "  +2.41%  "
$ws.Range("D38").Value = "'8.786"
$ws.Range("E38").Value = "  +5.10%  "
$ws.Range("D39").Value = "'5.090"
$ws.Range("E39").Value = "  +2.62%  "
$ws.Range("D40").Value = "'11.53"
$ws.Range("E40").Value = "  +1.37%  "
$ws.Range("D41").Value = "'0.6326"
$ws.Range("E41").Value = "  +2.42%  "
$ws.Range("E42").Value = "  +5.69%  "
$ws.Range("D43").Value = "'1.0000"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").Value = "'1.393"
$ws.Range("E44").Value = "  -1.88%  "
$ws.Range("D45").Value = "'13.53"
$ws.Range("E45").Value = "  +1.90%  "
$ws.Range("D46").Value = "'0.5937"
$ws.Range("E46").Value = "  +1.97%  "
$ws.Range("D47").Value = "'3.668"
$ws.Range("E47").Value = "  +0.92%  "
$ws.Range("D48").Value = "'125.80"
$ws.Range("E48").Value = "  +3.54%  "
$ws.Range("D49").Value = "'1.983"
$ws.Range("E49").Value = "  +3.87%  "
$ws.Range("D50").Value = "'1.170"
$ws.Range("E50").Value = "  +3.29%  "
$ws.Range("D51").Value = "'0.06937"
$ws.Range("E51").Value = "  +2.86%  "
